$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 200
$ws.Range("I29").Value = 200
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 600
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -319
$ws.Range("N29").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2014
$ws.Range("I40").Value = 1962.75
$ws.Range("J40").Value = 2059.5557
$ws.Range("K40").Value = 1962.75
$ws.Range("L40").Value = 2059.5557
$ws.Range("M40").Value = -1787.75
$ws.Range("N40").Value = -2409.5557

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 6530
$ws.Range("I113").Value = 3977.5
$ws.Range("K113").Value = 3977.5
$ws.Range("M113").Value = -723.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 8589.647000000001
$ws.Range("I74").Value = 1261
$ws.Range("J74").Value = 19059.143
$ws.Range("K74").Value = 1261
$ws.Range("L74").Value = 19059.143
$ws.Range("M74").Value = -387
$ws.Range("N74").Value = -20807.143

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 8589.647000000001
$ws.Range("I77").Value = 1261
$ws.Range("J77").Value = 19059.143
$ws.Range("K77").Value = 6305
$ws.Range("L77").Value = 95295.715
$ws.Range("M77").Value = -1937
$ws.Range("N77").Value = -104031.715

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3792.2593
$ws.Range("I122").Value = 3122.1765
$ws.Range("J122").Value = 4931.4
$ws.Range("K122").Value = 9366.529500000001
$ws.Range("L122").Value = 14794.2
$ws.Range("M122").Value = -6916.529500000001
$ws.Range("N122").Value = -19694.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4019.4
$ws.Range("I132").Value = 3590.611
$ws.Range("K132").Value = 10771.833
$ws.Range("M132").Value = -8241.832999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 450
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 450
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 450
$ws.Range("N22").Value = -796
$ws.Range("M22").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2333.3333
$ws.Range("I86").Value = 4000
$ws.Range("K86").Value = 4000
$ws.Range("M86").Value = -2877

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2333.3333
$ws.Range("I89").Value = 4000
$ws.Range("K89").Value = 20000
$ws.Range("M89").Value = -14384

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H106").Value = 29372.5
$ws.Range("J106").Value = 29372.5
$ws.Range("L106").Value = 29372.5
$ws.Range("N106").Value = -31896.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 25138.4
$ws.Range("I68").Value = 31000.5
$ws.Range("J68").Value = 1690
$ws.Range("K68").Value = 93001.5
$ws.Range("L68").Value = 5070
$ws.Range("M68").Value = -92190.5
$ws.Range("N68").Value = -6692

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 25138.4
$ws.Range("I71").Value = 31000.5
$ws.Range("J71").Value = 1690
$ws.Range("K71").Value = 279004.5
$ws.Range("L71").Value = 15210
$ws.Range("M71").Value = -274948.5
$ws.Range("N71").Value = -23322

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H112").Value = 66670340
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 66670340
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 200011020
$ws.Range("N112").Value = -200013236
$ws.Range("M112").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 7753423.5
$ws.Range("J131").Value = 8131605.5
$ws.Range("L131").Value = 24394816.5
$ws.Range("N131").Value = -24404896.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 13890007
$ws.Range("J132").Value = 23810640
$ws.Range("L132").Value = 214295760
$ws.Range("N132").Value = -214300820

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 5062.2856
$ws.Range("I102").Value = 7345.3335
$ws.Range("J102").Value = 3350
$ws.Range("K102").Value = 7345.3335
$ws.Range("L102").Value = 3350
$ws.Range("M102").Value = -5723.3335
$ws.Range("N102").Value = -6594

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2346
$ws.Range("I122").Value = 2325.75
$ws.Range("J122").Value = 2386.5
$ws.Range("K122").Value = 6977.25
$ws.Range("L122").Value = 7159.5
$ws.Range("M122").Value = -4527.25
$ws.Range("N122").Value = -12059.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2605.3333
$ws.Range("I126").Value = 2070.8572
$ws.Range("K126").Value = 6212.571599999999
$ws.Range("M126").Value = -3742.571599999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2960.36
$ws.Range("I7").Value = 2140.8
$ws.Range("K7").Value = 2140.8
$ws.Range("M7").Value = -2028.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 3000
$ws.Range("J18").Value = 3000
$ws.Range("L18").Value = 3000
$ws.Range("N18").Value = -3344

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2864.72
$ws.Range("I40").Value = 2013.2222
$ws.Range("J40").Value = 3343.6875
$ws.Range("K40").Value = 2013.2222
$ws.Range("L40").Value = 3343.6875
$ws.Range("M40").Value = -1877.2222
$ws.Range("N40").Value = -3615.6875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H70").Value = 19500
$ws.Range("J70").Value = 19500
$ws.Range("L70").Value = 19500
$ws.Range("N70").Value = -20040

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H73").Value = 19500
$ws.Range("J73").Value = 19500
$ws.Range("L73").Value = 19500
$ws.Range("N73").Value = -21372

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H106").Value = 16841.125
$ws.Range("J106").Value = 16841.125
$ws.Range("L106").Value = 16841.125
$ws.Range("N106").Value = -19365.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3864.1177
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 3864.1177
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 11592.3531
$ws.Range("N122").Value = -16492.3531
$ws.Range("M122").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2960.36
$ws.Range("I126").Value = 2140.8
$ws.Range("K126").Value = 6422.400000000001
$ws.Range("M126").Value = -3952.400000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1380.8
$ws.Range("I122").Value = 1380.8
$ws.Range("K122").Value = 4142.4
$ws.Range("M122").Value = -1692.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 42484.75
$ws.Range("I126").Value = 71838.5
$ws.Range("K126").Value = 215515.5
$ws.Range("M126").Value = -213045.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1322.6818
$ws.Range("I136").Value = 748.5
$ws.Range("K136").Value = 2245.5
$ws.Range("M136").Value = 304.5
